$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date values from 45184 to 45186 for rows 2-57
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# Add friendly display text (the "Beteckning" / column A value) as the second
# argument of the HYPERLINK() formulas in columns S,T,U,V,W,X,Y for rows 2-14

# Row 2: A 55967-2019
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 55967-2019.xlsx", "A 55967-2019")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 55967-2019.png", "A 55967-2019")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/knärot/A 55967-2019.png", "A 55967-2019")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 55967-2019.docx", "A 55967-2019")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 55967-2019.docx", "A 55967-2019")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 55967-2019.docx", "A 55967-2019")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 55967-2019.docx", "A 55967-2019")'

# Row 3: A 15601-2021
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 15601-2021.xlsx", "A 15601-2021")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 15601-2021.png", "A 15601-2021")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 15601-2021.docx", "A 15601-2021")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 15601-2021.docx", "A 15601-2021")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 15601-2021.docx", "A 15601-2021")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 15601-2021.docx", "A 15601-2021")'

# Row 4: A 33441-2021
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 33441-2021.xlsx", "A 33441-2021")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 33441-2021.png", "A 33441-2021")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 33441-2021.docx", "A 33441-2021")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 33441-2021.docx", "A 33441-2021")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 33441-2021.docx", "A 33441-2021")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 33441-2021.docx", "A 33441-2021")'

# Row 5: A 25724-2023
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 25724-2023.xlsx", "A 25724-2023")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 25724-2023.png", "A 25724-2023")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 25724-2023.docx", "A 25724-2023")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 25724-2023.docx", "A 25724-2023")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 25724-2023.docx", "A 25724-2023")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 25724-2023.docx", "A 25724-2023")'

# Row 6: A 46277-2018
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 46277-2018.xlsx", "A 46277-2018")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 46277-2018.png", "A 46277-2018")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 46277-2018.docx", "A 46277-2018")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 46277-2018.docx", "A 46277-2018")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 46277-2018.docx", "A 46277-2018")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 46277-2018.docx", "A 46277-2018")'

# Row 7: A 62133-2019
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 62133-2019.xlsx", "A 62133-2019")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 62133-2019.png", "A 62133-2019")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 62133-2019.docx", "A 62133-2019")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 62133-2019.docx", "A 62133-2019")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 62133-2019.docx", "A 62133-2019")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 62133-2019.docx", "A 62133-2019")'

# Row 8: A 50415-2020
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 50415-2020.xlsx", "A 50415-2020")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 50415-2020.png", "A 50415-2020")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 50415-2020.docx", "A 50415-2020")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 50415-2020.docx", "A 50415-2020")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 50415-2020.docx", "A 50415-2020")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 50415-2020.docx", "A 50415-2020")'

# Row 9: A 22025-2021
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 22025-2021.xlsx", "A 22025-2021")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 22025-2021.png", "A 22025-2021")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 22025-2021.docx", "A 22025-2021")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 22025-2021.docx", "A 22025-2021")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 22025-2021.docx", "A 22025-2021")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 22025-2021.docx", "A 22025-2021")'

# Row 10: A 6832-2023
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 6832-2023.xlsx", "A 6832-2023")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 6832-2023.png", "A 6832-2023")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 6832-2023.docx", "A 6832-2023")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 6832-2023.docx", "A 6832-2023")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 6832-2023.docx", "A 6832-2023")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 6832-2023.docx", "A 6832-2023")'

# Row 11: A 61147-2019
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 61147-2019.xlsx", "A 61147-2019")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 61147-2019.png", "A 61147-2019")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 61147-2019.docx", "A 61147-2019")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 61147-2019.docx", "A 61147-2019")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 61147-2019.docx", "A 61147-2019")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 61147-2019.docx", "A 61147-2019")'

# Row 12: A 61354-2019
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 61354-2019.xlsx", "A 61354-2019")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 61354-2019.png", "A 61354-2019")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 61354-2019.docx", "A 61354-2019")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 61354-2019.docx", "A 61354-2019")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 61354-2019.docx", "A 61354-2019")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 61354-2019.docx", "A 61354-2019")'

# Row 13: A 42895-2020
$ws.Range("S13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 42895-2020.xlsx", "A 42895-2020")'
$ws.Range("T13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 42895-2020.png", "A 42895-2020")'
$ws.Range("V13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 42895-2020.docx", "A 42895-2020")'
$ws.Range("W13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 42895-2020.docx", "A 42895-2020")'
$ws.Range("X13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 42895-2020.docx", "A 42895-2020")'
$ws.Range("Y13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 42895-2020.docx", "A 42895-2020")'

# Row 14: A 66277-2020
$ws.Range("S14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 66277-2020.xlsx", "A 66277-2020")'
$ws.Range("T14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 66277-2020.png", "A 66277-2020")'
$ws.Range("V14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 66277-2020.docx", "A 66277-2020")'
$ws.Range("W14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 66277-2020.docx", "A 66277-2020")'
$ws.Range("X14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 66277-2020.docx", "A 66277-2020")'
$ws.Range("Y14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 66277-2020.docx", "A 66277-2020")'
